$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: copy H1's format (bold font, border, centered alignment) onto I1/J1,
# then set their text values.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2 through 35: I column = 1 (constant), J column = copy of H column value
for ($r = 2; $r -le 35; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
